$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualiza notas dos alunos - preenche a coluna C3 (D) para os 3 alunos
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 2

# Recalcula a planilha para que as somas (coluna G) reflitam os novos valores
$excel.Calculate()

# Atualiza a seleção ativa para D3, conforme ficou ao final da edição
$ws.Range("D3").Select()
